# Doing Updates for Financials
# A new fiscal-year column (FY2018, period ending 2018-12-31) is inserted
# immediately before column D on the single worksheet "FFBC". Every
# existing value in columns D:K shifts right to E:L, and the new column D
# is populated with the latest-year figures. Two historical cells
# (previously-reported cumulative figures in the Cash Flow section) are
# also restated as part of the same update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FFBC")

# 1) Insert a new column before D; this shifts D:K -> E:L (values, formulas
#    and styles all move together).
$ws.Range("D1").EntireColumn.Insert()

# 2) The freshly inserted column D has no formatting of its own (Excel
#    copies it from the column to the left, i.e. column C). Copy the
#    number formats from column E (which now holds what used to be in D)
#    back onto D so the new column matches the rest of the data columns.
$fmtSource = $ws.Range("E7:E102")
$fmtDest = $ws.Range("D7:D102")
$fmtSource.Copy()
$fmtDest.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# 3) Populate the new column D with the FY2018 figures.
$newColumnValues = @{
    7  = 43465
    8  = 540400
    13 = 0
    14 = 0
    15 = 0
    17 = 105700
    18 = 434600
    20 = -220400
    21 = 238400
    22 = 0
    23 = 214200
    24 = 41600
    25 = 0
    26 = 172600
    27 = 172600
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 220400
    33 = 172600
    34 = 0
    35 = 172600
    38 = 43465
    41 = 236200
    42 = 153400
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 215700
    49 = 921100
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 13986700
    57 = 142000
    58 = 0
    59 = 0
    60 = 0
    61 = 570700
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 11908400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 600000
    73 = 0
    74 = 0
    75 = 0
    76 = 2078200
    77 = 0
    80 = 43465
    81 = 172600
    83 = 24200
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 261400
    91 = -18200
    92 = 0
    93 = 0
    94 = -205800
    96 = -79700
    97 = 0
    98 = 0
    99 = 0
    100 = 30000
    101 = 0
    102 = 85600
}

foreach ($row in $newColumnValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $newColumnValues[$row]
}

# Rows where the figure is reported as "NA" (text) rather than a number,
# same as the rest of the row.
$naRows = @(9, 10, 12)
foreach ($row in $naRows) {
    $ws.Cells.Item($row, 4).Value = "NA"
}

# 4) Two cash-flow rows were restated at the same time: the prior two
#    periods' figures changed (not just shifted) alongside the new column.
$ws.Cells.Item(89, 5).Value = 123500
$ws.Cells.Item(89, 6).Value = 142600

$ws.Cells.Item(94, 5).Value = -429500
$ws.Cells.Item(94, 6).Value = -313300
